# Sync attendance_reports: normalize "Recorded By" (column G) ordering so that
# the last contributor listed is rotated to the front of the comma-separated
# list (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System").
#
# The single exception is the literal value "admin@admin.com, System", which
# is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G ("Recorded By") is column index 7.
$col = 7

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    $text = [string]$val

    if ($text -eq "admin@admin.com, System") {
        continue
    }

    $parts = $text -split ", "

    if ($parts.Length -le 1) {
        continue
    }

    $lastPart = $parts[$parts.Length - 1]
    $remaining = $parts[0..($parts.Length - 2)]
    $newParts = @($lastPart) + $remaining
    $newText = [string]::Join(", ", $newParts)

    $cell.Value = $newText
}
